# Fruta / hortaliza, semanal
# Insert a new weekly data row at row 5 (pushing existing rows 5-96 down to 6-97)
# and populate it with the new record's values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 5; this shifts rows 5..96 down to 6..97 and
# carries over the existing formatting (e.g. the date number format on column D).
$ws.Rows("5:5").Insert()

# Columns that stay constant across every record in this sheet.
$ws.Range("A5").Value = 7
$ws.Range("B5").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C5").Value = "Ñuble"
$ws.Range("E5").Value = 16
$ws.Range("F5").Value = 100112031
$ws.Range("G5").Value = "Poroto verde"
$ws.Range("I5").Value = "Primera"
$ws.Range("Q5").Value = 25
$ws.Range("R5").Value = "Hortaliza"

# New record's specific values.
$ws.Range("D5").Value = 44812
$ws.Range("H5").Value = "Magnum"
$ws.Range("J5").Value = 60
$ws.Range("K5").Value = 30000
$ws.Range("L5").Value = 30000
$ws.Range("M5").Value = 30000
$ws.Range("N5").Value = "$/malla 25 kilos"
$ws.Range("O5").Value = "Perú"
$ws.Range("P5").Value = 1200
